$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values that look like plain numbers stay as text,
# matching the source data which always stores Price/Volume as strings.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.816.82"
$ws.Range("E2").Value = "  +1.27%  "
$ws.Range("D3").Value = "1.885.76"
$ws.Range("E3").Value = "  +1.10%  "
$ws.Range("D4").Value = "1.018"
$ws.Range("E4").Value = "  +1.62%  "
$ws.Range("D5").Value = "333.84"
$ws.Range("E5").Value = "  +1.23%  "
$ws.Range("D6").Value = "1.014"
$ws.Range("E6").Value = "  +1.31%  "
$ws.Range("D7").Value = "0.4698"
$ws.Range("E7").Value = "  -0.38%  "
$ws.Range("D8").Value = "0.3912"
$ws.Range("E8").Value = "  -1.61%  "
$ws.Range("D9").Value = "47.86"
$ws.Range("E9").Value = "  +0.37%  "
$ws.Range("D10").Value = "0.08057"
$ws.Range("E10").Value = "  +0.33%  "
$ws.Range("D11").Value = "1.015"
$ws.Range("E11").Value = "  -0.40%  "
$ws.Range("D12").Value = "21.92"
$ws.Range("E12").Value = "  +1.37%  "
$ws.Range("D13").Value = "1.898.70"
$ws.Range("E13").Value = "  +1.58%  "
$ws.Range("D14").Value = "5.951"
$ws.Range("E14").Value = "  -0.02%  "
$ws.Range("D15").Value = "7.060"
$ws.Range("E15").Value = "  -1.87%  "
$ws.Range("D16").Value = "1.018"
$ws.Range("E16").Value = "  +1.69%  "
$ws.Range("D17").Value = "0.06747"
$ws.Range("E17").Value = "  +2.75%  "
$ws.Range("D18").Value = "87.22"
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("D19").Value = "0.00001045"
$ws.Range("E19").Value = "  +0.47%  "
$ws.Range("D20").Value = "17.17"
$ws.Range("E20").Value = "  -1.14%  "
$ws.Range("D21").Value = "1.013"
$ws.Range("E21").Value = "  +1.18%  "
$ws.Range("D22").Value = "27.864.28"
$ws.Range("E22").Value = "  +1.41%  "
$ws.Range("D23").Value = "5.494"
$ws.Range("E23").Value = "  -0.42%  "
$ws.Range("D24").Value = "10.92"
$ws.Range("E24").Value = "  -0.49%  "
$ws.Range("D25").Value = "2.342"
$ws.Range("E25").Value = "  +1.84%  "
$ws.Range("D26").Value = "2.133.61"
$ws.Range("E26").Value = "  +2.45%  "
$ws.Range("D27").Value = "159.62"
$ws.Range("E27").Value = "  +3.65%  "
$ws.Range("D28").Value = "19.98"
$ws.Range("E28").Value = "  -1.58%  "
$ws.Range("D29").Value = "2.081"
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").Value = "5.495"
$ws.Range("E30").Value = "  -0.85%  "
$ws.Range("D31").Value = "121.66"
$ws.Range("E31").Value = "  -0.55%  "
$ws.Range("D32").Value = "0.9706"
$ws.Range("E32").Value = "  +1.44%  "
$ws.Range("D33").Value = "0.09459"
$ws.Range("E33").Value = "  -0.36%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "3.645"
$ws.Range("E34").Value = "  +1.43%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "1.412"
$ws.Range("E35").Value = "  -3.73%  "
$ws.Range("D36").Value = "5.335"
$ws.Range("E36").Value = "  +0.61%  "
$ws.Range("D37").Value = "0.06118"
$ws.Range("E37").Value = "  +0.35%  "
$ws.Range("D38").Value = "0.02255"
$ws.Range("E38").Value = "  +0.47%  "
$ws.Range("D39").Value = "1.211"
$ws.Range("E39").Value = "  -0.82%  "
$ws.Range("D40").Value = "0.5954"
$ws.Range("E40").Value = "  -0.47%  "
$ws.Range("D41").Value = "7.976"
$ws.Range("E41").Value = "  -1.47%  "
$ws.Range("D42").Value = "0.1885"
$ws.Range("E42").Value = "  -0.48%  "
$ws.Range("D43").Value = "10.24"
$ws.Range("E43").Value = "  -0.97%  "
$ws.Range("D44").Value = "1.264"
$ws.Range("E44").Value = "  +0.20%  "
$ws.Range("D45").Value = "0.5667"
$ws.Range("E45").Value = "  -0.46%  "
$ws.Range("D46").Value = "12.12"
$ws.Range("E46").Value = "  -0.53%  "
$ws.Range("D47").Value = "3.405"
$ws.Range("E47").Value = "  -0.31%  "
$ws.Range("D48").Value = "1.927"
$ws.Range("E48").Value = "  -0.75%  "
$ws.Range("D49").Value = "0.06906"
$ws.Range("E49").Value = "  +1.87%  "
$ws.Range("D50").Value = "113.64"
$ws.Range("E50").Value = "  +3.19%  "
$ws.Range("D51").Value = "1.069"
$ws.Range("E51").Value = "  +0.57%  "
